# Updates cryptos list data (Price / Volume(1h) columns) to match the
# latest scrape, per commit "Updated cryptos list on Sun Oct 13 12:51:25 UTC 2024".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.729.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.464.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.112"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.636.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.466.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.91"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.19"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +15.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "640.27"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.588.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0974"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -15.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.74"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.31"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0307"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -17.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "153.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.30"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.38"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("E51").Value = "  -1.29%  "
